$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 288525.12
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 64
$ws.Range("H64").Value = 5291.1665
$ws.Range("J64").Value = 5437
$ws.Range("L64").Value = 5437
$ws.Range("N64").Value = -5933

# Row 67
$ws.Range("H67").Value = 5291.1665
$ws.Range("J67").Value = 5437
$ws.Range("L67").Value = 5437
$ws.Range("N67").Value = -7153

# Row 70
$ws.Range("H70").Value = 3498.5
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 3664.6667
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 10994.0001
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -11534.0001

# Row 73
$ws.Range("H73").Value = 3498.5
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 3664.6667
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 10994.0001
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -12866.0001

# Row 98
$ws.Range("H98").Value = 672.1053000000001
$ws.Range("I98").Value = 505.55356
$ws.Range("K98").Value = 505.55356
$ws.Range("M98").Value = 992.4464399999999

# Row 106
$ws.Range("H106").Value = 12150.818
$ws.Range("I106").Value = 2239.4285
$ws.Range("K106").Value = 2239.4285
$ws.Range("M106").Value = -1608.4285

# Row 112
$ws.Range("H112").Value = 1350.6
$ws.Range("I112").Value = 1029.6666
$ws.Range("K112").Value = 3088.9998
$ws.Range("M112").Value = -1980.9998

# Row 121
$ws.Range("H121").Value = 4369.7334
$ws.Range("J121").Value = 4369.7334
$ws.Range("L121").Value = 13109.2002
$ws.Range("N121").Value = -16603.2002

# Row 122
$ws.Range("H122").Value = 672.1053000000001
$ws.Range("I122").Value = 505.55356
$ws.Range("K122").Value = 1516.66068
$ws.Range("M122").Value = 933.33932

# Row 132
$ws.Range("H132").Value = 1384.8684
$ws.Range("I132").Value = 1404.0834
$ws.Range("K132").Value = 4212.2502
$ws.Range("M132").Value = -1682.2502

# Row 137
$ws.Range("H137").Value = 2151.5
$ws.Range("I137").Value = 2263.3845
$ws.Range("J137").Value = 1666.6666
$ws.Range("K137").Value = 6790.1535
$ws.Range("L137").Value = 4999.9998
$ws.Range("M137").Value = -4240.1535
$ws.Range("N137").Value = -10099.9998

# Row 138
$ws.Range("H138").Value = 1581.5205
$ws.Range("I138").Value = 1359.9565
$ws.Range("J138").Value = 1959
$ws.Range("K138").Value = 4079.8695
$ws.Range("L138").Value = 5877
$ws.Range("M138").Value = 1060.1305
$ws.Range("N138").Value = -16157

# Row 141
$ws.Range("H141").Value = 18489.16
$ws.Range("I141").Value = 19992.72
$ws.Range("K141").Value = 59978.16
$ws.Range("M141").Value = -54798.16

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1706.6296
$ws.Range("I74").Value = 1518.238
$ws.Range("J74").Value = 2366
$ws.Range("K74").Value = 1518.238
$ws.Range("L74").Value = 2366
$ws.Range("M74").Value = -644.2380000000001
$ws.Range("N74").Value = -4114

# Row 77
$ws.Range("H77").Value = 1706.6296
$ws.Range("I77").Value = 1518.238
$ws.Range("J77").Value = 2366
$ws.Range("K77").Value = 7591.190000000001
$ws.Range("L77").Value = 11830
$ws.Range("M77").Value = -3223.190000000001
$ws.Range("N77").Value = -20566

# Row 88
$ws.Range("H88").Value = 835.28
$ws.Range("I88").Value = 538.0833
$ws.Range("J88").Value = 1109.6154
$ws.Range("K88").Value = 538.0833
$ws.Range("L88").Value = 1109.6154
$ws.Range("M88").Value = -132.0833
$ws.Range("N88").Value = -1921.6154

# Row 91
$ws.Range("H91").Value = 835.28
$ws.Range("I91").Value = 538.0833
$ws.Range("J91").Value = 1109.6154
$ws.Range("K91").Value = 538.0833
$ws.Range("L91").Value = 1109.6154
$ws.Range("M91").Value = 865.9167
$ws.Range("N91").Value = -3917.6154

# Row 132
$ws.Range("H132").Value = 1546.931
$ws.Range("I132").Value = 1462.3928
$ws.Range("K132").Value = 4387.178400000001
$ws.Range("M132").Value = -1857.178400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 16672.234
$ws.Range("I20").Value = 17245.3
$ws.Range("J20").Value = 15853.571
$ws.Range("K20").Value = 17245.3
$ws.Range("L20").Value = 15853.571
$ws.Range("M20").Value = -16998.3
$ws.Range("N20").Value = -16347.571

# Row 107
$ws.Range("H107").Value = 25915.62
$ws.Range("I107").Value = 37372.43
$ws.Range("K107").Value = 37372.43
$ws.Range("M107").Value = -35452.43

# Row 134
$ws.Range("H134").Value = 1116.8889
$ws.Range("I134").Value = 1116.8889
$ws.Range("K134").Value = 3350.6667
$ws.Range("M134").Value = -815.6666999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 86969.31
$ws.Range("I31").Value = 147213.58
$ws.Range("J31").Value = 40112.668
$ws.Range("K31").Value = 147213.58
$ws.Range("L31").Value = 40112.668
$ws.Range("M31").Value = -146918.58
$ws.Range("N31").Value = -40702.668

# Row 34
$ws.Range("H34").Value = 86969.31
$ws.Range("I34").Value = 147213.58
$ws.Range("J34").Value = 40112.668
$ws.Range("K34").Value = 147213.58
$ws.Range("L34").Value = 40112.668
$ws.Range("M34").Value = -147011.58
$ws.Range("N34").Value = -40516.668

# Row 62
$ws.Range("H62").Value = 5356.2856
$ws.Range("I62").Value = 5499
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 5499
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -4875
$ws.Range("N62").Value = -5748

# Row 65
$ws.Range("H65").Value = 5356.2856
$ws.Range("I65").Value = 5499
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 27495
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -24375
$ws.Range("N65").Value = -28740

# Row 68
$ws.Range("H68").Value = 24998.572
$ws.Range("J68").Value = 24998.572
$ws.Range("L68").Value = 24998.572
$ws.Range("N68").Value = -26496.572

# Row 71
$ws.Range("H71").Value = 24998.572
$ws.Range("J71").Value = 24998.572
$ws.Range("L71").Value = 74995.716
$ws.Range("N71").Value = -82483.716

# Row 134
$ws.Range("H134").Value = 2574.0857
$ws.Range("I134").Value = 2406.8
$ws.Range("K134").Value = 7220.400000000001
$ws.Range("M134").Value = -4685.400000000001

# Row 141
$ws.Range("H141").Value = 137581.86
$ws.Range("J141").Value = 148012.67
$ws.Range("L141").Value = 148012.67
$ws.Range("N141").Value = -158372.67

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 754.3077
$ws.Range("J107").Value = 650.6667
$ws.Range("L107").Value = 1952.0001
$ws.Range("N107").Value = -5792.0001

# Row 132
$ws.Range("H132").Value = 2356.762
$ws.Range("I132").Value = 1431.8889
$ws.Range("K132").Value = 12887.0001
$ws.Range("M132").Value = -10357.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 128
$ws.Range("H128").Value = 38333
$ws.Range("J128").Value = 38333
$ws.Range("L128").Value = 38333
$ws.Range("N128").Value = -48293

# Row 132
$ws.Range("H132").Value = 2725.353
$ws.Range("I132").Value = 2166.7144
$ws.Range("K132").Value = 6500.1432
$ws.Range("M132").Value = -3970.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6005.9346
$ws.Range("I7").Value = 6030.657
$ws.Range("K7").Value = 6030.657
$ws.Range("M7").Value = -5918.657

# Row 122
$ws.Range("H122").Value = 3678.6667
$ws.Range("I122").Value = 2700.2
$ws.Range("K122").Value = 8100.599999999999
$ws.Range("M122").Value = -5650.599999999999

# Row 126
$ws.Range("H126").Value = 6005.9346
$ws.Range("I126").Value = 6030.657
$ws.Range("K126").Value = 18091.971
$ws.Range("M126").Value = -15621.971

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 540.4737
$ws.Range("I113").Value = 421.53845
$ws.Range("K113").Value = 1264.61535
$ws.Range("M113").Value = 905.38465

# Row 122
$ws.Range("H122").Value = 1644.1666
$ws.Range("I122").Value = 1664.1666
$ws.Range("K122").Value = 4992.4998
$ws.Range("M122").Value = -2542.4998

# Row 126
$ws.Range("H126").Value = 2919.36
$ws.Range("I126").Value = 2973.75
$ws.Range("J126").Value = 2701.8
$ws.Range("K126").Value = 8921.25
$ws.Range("L126").Value = 8105.400000000001
$ws.Range("M126").Value = -6451.25
$ws.Range("N126").Value = -13045.4
